$d = $word.ActiveDocument

# Locate the paragraph to remove entirely: "Improve string conversion APIs, ..."
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Improve string conversion APIs*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph to remove."
}

$target = $d.Paragraphs.Item($targetIndex)

# The paragraph immediately before it is where the "last edit" (_GoBack)
# bookmark needs to end up, right after its final run of text.
$prev = $d.Paragraphs.Item($targetIndex - 1)

# Remove the existing _GoBack bookmark (currently at the end of the x64
# trampoline paragraph) so it can be re-created at the new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

# Build a collapsed range at the very end of $prev's text (just before its
# paragraph mark) and use a temporary marker run to anchor a bookmark
# there, then shrink the marker back out, leaving a zero-length
# "_GoBack" bookmark at that exact position.
$anchor = $prev.Range.Duplicate
$anchor.MoveEnd(1, -1) | Out-Null
$anchor.Collapse(0) | Out-Null
$anchor.InsertAfter("IronGoBackMarker") | Out-Null

$marker = $d.Content
$marker.Find.Execute("IronGoBackMarker", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $marker) | Out-Null
$marker.Text = ""

# Finally, delete the whole paragraph (including its paragraph mark).
$target.Range.Delete() | Out-Null
